# ---------------------------------------------------------------------------
# "Correção de todos os artefatos"
#
# 1) Slide 5, shape 8 ("Google Shape;125;p17"): fix the bullet text
#       "Entrega de pedido(s)"  ->  "Entregar pedido"
#
# 2) Theme colors: the deck's live theme (the one used by the slide master /
#    all slides) gets the "Default" color palette (it previously carried the
#    "Simple Light" palette, which the edit moves elsewhere). We reproduce
#    the 12 theme colors via the slide's ThemeColorScheme, which edits the
#    shared theme part used by every slide (not a per-slide override).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Text fix -----------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shp = $s5.Shapes.Item(8)
$shp.TextFrame.TextRange.Runs(1).Text = "Entregar pedido"

# --- 2) Theme color fix ------------------------------------------------------
# RGBColor.RGB values use the OLE COLORREF byte order 0x00BBGGRR, i.e. the
# bytes of the hex RRGGBB color are reversed.
function BGR([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target ("Default") palette, in ThemeColorScheme order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$targetHex = @(
    "000000", "FFFFFF", "158158", "F3F3F3",
    "058DC7", "50B432", "ED561B", "EDEF00",
    "24CBE5", "64E572", "2200CC", "551A8B"
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = BGR $targetHex[$i - 1]
}
